$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "My Series" to "Data"
$ws.Name = "Data"

# Update the header cell text
$ws.Range("K1").Value = "Function Information"
